$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Correct the unit price of the 10uF capacitor (row 5, C1/C3/C4/C5/C6)
$ws.Range("E5").Value = 0.074

# Recalculate dependent formulas (F5 shared formula, F18 sum)
$excel.Calculate()

# Update the selected/active cell to reflect where the edit was made
$ws.Range("E5").Select()
